# "added Date type to import functionality"
# Adds a new column H to Sheet1: a Date-formatted value in H1 and a
# single-space text value (matching the existing shared string " ") in H2:H7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H1: a Date value (2015-01-01, Excel serial 42005) with a date number format
$ws.Range("H1").Value = 42005
$ws.Range("H1").NumberFormat = "mm-dd-yy"

# H2:H7: a literal single-space string value (reuses the existing shared string)
$ws.Range("H2").Value = " "
$ws.Range("H3").Value = " "
$ws.Range("H4").Value = " "
$ws.Range("H5").Value = " "
$ws.Range("H6").Value = " "
$ws.Range("H7").Value = " "

# Move the active selection to the newly-added cell H7
$ws.Range("H7").Select()
